$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D16").Value = "2016-03-09 05:00:13"
$wsZhCn.Range("G16").Value = "2016-03-09 05:01:14"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D16").Value = "2016-03-09 05:00:16"
$wsDeDe.Range("G16").Value = "2016-03-09 05:01:18"
